# Modificar registro de planificaciones y consulta de planificaciones
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekday columns (Mon-Fri) repeated across the 4 full weeks + partial last week,
# skipping the weekend columns (Sat/Sun) - matches the header row pattern.
$cols = @("B","C","D","E","F","I","J","K","L","M","P","Q","R","S","T","W","X","Y","Z","AA","AD","AE")

# Rows 3-12: fill every weekday cell with the default planning value "HR01"
foreach ($row in 3..12) {
    foreach ($col in $cols) {
        $ws.Range($col + $row).Value = "HR01"
    }
}

# A handful of cells hold a different (combined) planning value instead of the
# default - applied in the same order the author originally typed them in, so
# new shared-string entries land in the same order.
$ws.Range("B10").Value = "HR01, HR02"
$ws.Range("D9").Value  = "HR01, C4"
$ws.Range("M6").Value  = "FERIADO 1"
$ws.Range("M9").Value  = "HR01, FERIADO 1"
$ws.Range("C6").Value  = "HR01,LHORARIO2"

# Move the current selection, as left by the author after editing
$ws.Range("E29").Select() | Out-Null
